# Methods and keywords edits
# The two "standalone" keyword rows (Phytoplankton growth / Zooplankton growth)
# that had no keywordThesaurus value are removed from the Keywords sheet.
# Deleting the full rows shifts the remaining "Northeast U.S. Continental Shelf" /
# "NOAA Large Marine Ecosystems" row up, and Excel automatically compacts the
# shared-strings table, which is what drives all of the value-index shifts
# seen elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

# Make Keywords the active sheet (matches activeTab moving to it / tabSelected
# moving off ColumnHeaders).
$ws.Activate() | Out-Null

# Remove the two rows that only held a bare keyword with no thesaurus.
$ws.Rows("12:13").Delete() | Out-Null

# Restore the expected selection / scroll position on the now-shorter sheet.
$ws.Range("A12:XFD13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
